$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness column (C) values for the log data.
# C2 (Generation 0) gets a distinct value, C3:C252 (Generations 1-250) share another value.
$ws.Range("C2").Value = 7581
$ws.Range("C3:C252").Value = 7573
